# Insert a new data row at row 227 (pushing the existing rows 227-298 down
# to 228-299), and populate it with the new weekly price-report entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("227:227").Insert()

$ws.Range("A227").Value = 10
$ws.Range("B227").Value = "Vega Modelo de Temuco"
$ws.Range("C227").Value = "La Araucanía"
$ws.Range("D227").Value = 45093
$ws.Range("E227").Value = 9
$ws.Range("F227").Value = 100112013
$ws.Range("G227").Value = "Alcachofa"
$ws.Range("H227").Value = "Española"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 600
$ws.Range("K227").Value = 550
$ws.Range("L227").Value = 550
$ws.Range("M227").Value = 550
$ws.Range("N227").Value = "`$/unidad"
$ws.Range("O227").Value = "Provincia de Limarí"
$ws.Range("P227").Value = 550
$ws.Range("Q227").Value = 1
$ws.Range("R227").Value = "Hortaliza"
